$d = $word.ActiveDocument

$replacements = @(
    @("e||----", "e||-"),
    @("B||----", "B||-"),
    @("G||----", "G||-"),
    @("D||\---", "D||-"),
    @("A||----", "A||-"),
    @("E||-/--", "E||-")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
